# Replace sleep with safer wait
#
# Three code-sample text boxes each contain a paragraph whose text "sleep 5"
# was typed as three separate runs ("s" / "leep " / "5"). The author retyped
# that paragraph as "wait": deleting the leading "s" first drops the first
# run (and its rPr) entirely, then overwriting the remaining "leep 5" with
# "wait" collapses the rest into a single run using the formatting of the
# run that used to read "leep ".

$p = $ppt.ActivePresentation

# --- Slide 10 -- shape "Rectangle 4" (lang="nl-BE") ---
$slide10 = $p.Slides.Item(10)
$shape10 = $slide10.Shapes.Item("Rectangle 4")
$tr10 = $shape10.TextFrame.TextRange
$full10 = $tr10.Text
$idx10 = $full10.IndexOf("sleep 5")
$start10 = $tr10.Start + $idx10
$lead10 = $tr10.Characters($start10, 1)
$lead10.Text = ""
$rest10 = $tr10.Characters($start10, 6)
$rest10.Text = "wait"

# --- Slide 11 -- shape "Rectangle 2", last paragraph (lang="en-US") ---
$slide11 = $p.Slides.Item(11)
$shape11 = $slide11.Shapes.Item("Rectangle 2")
$tr11 = $shape11.TextFrame.TextRange
$lastPara11 = $tr11.Paragraphs($tr11.Paragraphs().Count)
$full11 = $lastPara11.Text
$idx11 = $full11.IndexOf("sleep 5")
$start11 = $lastPara11.Start + $idx11
$lead11 = $tr11.Characters($start11, 1)
$lead11.Text = ""
$rest11 = $tr11.Characters($start11, 6)
$rest11.Text = "wait"

# --- Slide 19 -- shape "Rectangle 3", last paragraph (lang="en-US") ---
$slide19 = $p.Slides.Item(19)
$shape19 = $slide19.Shapes.Item("Rectangle 3")
$tr19 = $shape19.TextFrame.TextRange
$lastPara19 = $tr19.Paragraphs($tr19.Paragraphs().Count)
$full19 = $lastPara19.Text
$idx19 = $full19.IndexOf("sleep 5")
$start19 = $lastPara19.Start + $idx19
$lead19 = $tr19.Characters($start19, 1)
$lead19.Text = ""
$rest19 = $tr19.Characters($start19, 6)
$rest19.Text = "wait"
